$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Card effect text updates (sharedStrings content) ---
# Row 2-6 keep their original text (only reordering happens automatically under
# the hood because other strings were removed/added), so we don't need to touch
# them explicitly. Only rows whose effect text actually changed need updating.

# D9 - 商人 (Shopkeeper) "Trade machine": every option now optional, pay-by-discard option
$ws.Range("D9").Value = "多选，可重复：①支付1金币，获得1道具点。②弃置1张战利品牌，获得1金币。③支付3金币，从遗物牌堆翻开3张牌，选其中1张获得。"

# D10 - 训练师 (Trainer): cost reduced from 5 time units to 3
$ws.Range("D10").Value = "消耗3时间，然后使用1任意属性或将1张手牌洗回主牌堆：获得1技能点。使用1张《智力》发动本牌时，可以少消耗2时间。"

# D8 - 冒险者尸体 (Dead adventurer): now also costs an attribute/hand card
$ws.Range("D8").Value = "消耗2时间，然后使用1任意属性或将1张手牌洗回主牌堆：将主牌堆第1张怪物牌放在房间区任意非空列顶端，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆顶的1张遗物牌。"

# D11 - 祭坛 (Altar): extra alternative reward added
$ws.Range("D11").Value = "将1张手牌送墓，或弃置1张战利品牌，或受到1点伤害：获得1道具点，或获得遭遇牌堆第1张战利品牌。"

# D7 - 宝箱 (Chest): reveal 5 cards instead of 3
$ws.Range("D7").Value = "消耗1时间，或使用1任意属性，或将1张手牌洗回主牌堆：从遭遇牌堆翻开5张牌，获得其中的战利品牌。使用1张《敏捷》发动本牌时，可以额外翻开2张牌。"

# --- Selection state ---
$ws.Range("D7").Select() | Out-Null
